$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph references can go stale once the document is mutated (e.g. by
# Paste), so every lookup below re-scans $d.Paragraphs fresh by index
# instead of holding on to an object across a mutating call.
# ---------------------------------------------------------------------------

$loopIndex = 0
$sourceIndex = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($loopIndex -eq 0 -and $t.TrimEnd([char]13, [char]7) -eq "If not retirement, loop") {
        $loopIndex = $i
    }
    if ($sourceIndex -eq 0 -and $t.StartsWith("Make class")) {
        $sourceIndex = $i
    }
    $i = $i + 1
}

# ---------------------------------------------------------------------------
# 1) Insert a new ListParagraph bullet ("Fix code so that if comment seen in
#    'board.txt', ignore.") at ilvl=1 / numId=1, right before the
#    "If not retirement, loop" item.
#
#    The interop layer has no direct way to reassign a paragraph to a
#    different numId, so instead we copy an existing paragraph that already
#    carries the numId=1/ilvl=1 numbering (the "Make class 'Retirement'"
#    bullet) via the clipboard, paste a clone of it in place, then
#    overwrite its text.
# ---------------------------------------------------------------------------
$sourcePara = $d.Paragraphs.Item($sourceIndex)
$sourcePara.Range.Copy()

$loopPara = $d.Paragraphs.Item($loopIndex)
$insertionPoint = $d.Range($loopPara.Range.Start, $loopPara.Range.Start)
$insertionPoint.Paste()

# The pasted paragraph now occupies $loopIndex, pushing "If not retirement,
# loop" one slot later.
$newPara = $d.Paragraphs.Item($loopIndex)
$q1 = [char]0x2018
$q2 = [char]0x2019
$newPara.Range.Text = "Fix code so that if comment seen in " + $q1 + "board.txt" + $q2 + ", ignore."

# ---------------------------------------------------------------------------
# 2) The "What to do when land on ... spaces-make methods" paragraph is
#    immediately followed by a "Test Spin to Win" paragraph that carries the
#    _GoBack bookmark. Merge the two paragraphs (delete the paragraph mark
#    between them), drop the "Test Spin to Win" text, and collapse the
#    run text down to the simpler wording, keeping the bookmark intact.
# ---------------------------------------------------------------------------
$landIndex = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($landIndex -eq 0 -and $p.Range.Text.StartsWith("What to do when land on")) {
        $landIndex = $i
    }
    $i = $i + 1
}

$landPara = $d.Paragraphs.Item($landIndex)
$markStart = $landPara.Range.End - 1
$d.Range($markStart, $markStart + 1).Delete()

$landPara = $d.Paragraphs.Item($landIndex)
$landPara.Range.Find.Execute("Test Spin to Win", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$landPara = $d.Paragraphs.Item($landIndex)
$landPara.Range.Find.Execute(" House and Retirement", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Delete the trailing to-do items that got removed in this revision:
#    "Add career object to player", "Fix up spaces", "Change career",
#    "Family Stop space", and the blank paragraph that used to close the
#    document body. Everything from right after the merged paragraph
#    through the end of the document goes away.
# ---------------------------------------------------------------------------
$landPara = $d.Paragraphs.Item($landIndex)
$tailStart = $landPara.Range.End
$d.Range($tailStart, $d.Content.End).Delete()
